$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the week date headers for the first two week columns
$ws.Range("B1").Value = "1/2 - 1/6"
$ws.Range("C1").Value = "1/9 - 1/13"

# Update recorded hours for each team member (rows 3-6) for the two weeks
# Joeseph Sedutto
$ws.Range("B3").Value = 4
$ws.Range("C3").Value = 1.25

# Matthew Handley
$ws.Range("B4").Value = 4
$ws.Range("C4").Value = 1.25

# Austin Cardosi
$ws.Range("B5").Value = 4

# Kathryn Swineford
$ws.Range("B6").Value = 4
$ws.Range("C6").Value = 1.25

# Update the selected cell to match the saved view state
$ws.Range("D11").Select()
